# Update GSC export data for main domain (HTTPS.xlsx) — append 3 new
# daily rows (2025-11-15, 2025-11-16, 2025-11-17) to the "Chart" sheet.
#
# The date values must land as literal text (matching the existing rows,
# which are stored as shared strings), not as auto-converted date
# serials, so each date cell is briefly switched to a text number format
# before the write and then cleared back to the default format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-DateRow($row, $dateText, $nonHttps, $https) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dateText
    $cell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $nonHttps
    $ws.Cells.Item($row, 3).Value = $https
}

Set-DateRow 41 "2025-11-15" 0 37
Set-DateRow 42 "2025-11-16" 0 35
Set-DateRow 43 "2025-11-17" 0 30
